$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: ADA LUZ UTRIA NAVARRO record (moved to top)
$ws.Range("C16").Value = "1048943696"
$ws.Range("D16").Value = "ADA LUZ UTRIA NAVARRO"
$ws.Range("E16").Value = "1808"
$ws.Range("F16").Value = 16666
$ws.Range("G16").Value = 1000000

# Row 17: SANDRA MILENA CASTILLO ORTIZ, period 2210
$ws.Range("C17").Value = "1143343913"
$ws.Range("D17").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E17").Value = "2210"
$ws.Range("F17").Value = 9333
$ws.Range("G17").Value = 1423500

# Row 18: SANDRA MILENA CASTILLO ORTIZ, period 2211
$ws.Range("C18").Value = "1143343913"
$ws.Range("D18").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E18").Value = "2211"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1423500

# Row 19: SANDRA MILENA CASTILLO ORTIZ, period 2212
$ws.Range("C19").Value = "1143343913"
$ws.Range("D19").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E19").Value = "2212"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1423500
